# Update Cleveland projections row 4 (data row) with new values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 6814
$ws.Range("B4").Value = 1073
$ws.Range("C4").Value = 27201
$ws.Range("D4").Value = 1.169502288329519
$ws.Range("E4").Value = 0.1842391304347826
$ws.Range("F4").Value = 4.668432054919907
$ws.Range("G4").Value = 0.09815011441647593
$ws.Range("H4").Value = 0.02351258581235698
$ws.Range("I4").Value = 0.1945364830091533
$ws.Range("J4").Value = 0.0688388947368421
$ws.Range("K4").Value = 0.04706184210526316
$ws.Range("L4").Value = 0.1281685545263158
$ws.Range("M4").Value = 10170
$ws.Range("N4").Value = 608
$ws.Range("O4").Value = 99685
$ws.Range("P4").Value = 1.824931291866028
$ws.Range("Q4").Value = 0.1090909090909091
$ws.Range("R4").Value = 17.88620904936603
$ws.Range("S4").Value = 0.1105020861244019
$ws.Range("T4").Value = 0.01238038277511962
$ws.Range("U4").Value = 0.3384172753708135
$ws.Range("V4").Value = 0.09608930315789471
$ws.Range("W4").Value = 0.04900394736842105
$ws.Range("X4").Value = 0.3985519831968421
